$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary fields ---
$ws.Range("C2").Value = 45355
$ws.Range("C3").Value = "'2556845"
$ws.Range("C4").Value = "DISEÑO E INTEGRACION DE MULTIMEDIA"
$ws.Range("C6").Value = 44760
$ws.Range("C7").Value = 45124

# --- Update existing data rows (11-15) ---
$ws.Range("B11").Value = 1005178211
$ws.Range("C11").Value = "JUAN CAMILO"
$ws.Range("D11").Value = "DELGADO CARRASCAL"
$ws.Range("E11").Value = "CERTIFICADO"

$ws.Range("B12").Value = 1005181992
$ws.Range("C12").Value = "JOHAN"
$ws.Range("D12").Value = "VARGAS CALDERIN"
$ws.Range("E12").Value = "POR CERTIFICAR"

$ws.Range("B13").Value = 1005184329
$ws.Range("C13").Value = "WILLIAM ANDRES"
$ws.Range("D13").Value = "LOPEZ RIOS"
$ws.Range("E13").Value = "CANCELADO"

$ws.Range("B14").Value = 1005185919
$ws.Range("C14").Value = "SEBASTIAN"
$ws.Range("D14").Value = "PERTUZ SAMPAYO"
$ws.Range("E14").Value = "CERTIFICADO"

$ws.Range("B15").Value = 1005220651
$ws.Range("C15").Value = "BRAYAN EDUARDO"
$ws.Range("D15").Value = "BADILLO HERRERA"
$ws.Range("E15").Value = "CERTIFICADO"

# --- New rows 16-27 ---
$newRows = @(
    @("CC", 1005239745, "SARAY DUVIANA", "UNRIZA JAIMES", "CERTIFICADO"),
    @("CC", 1005241421, "CLARA LUCIA", "RUIZ MONSALVE", "RETIRO VOLUNTARIO"),
    @("CC", 1043962939, "DANNA KAROLAY", "RESTREPO SOSA", "CERTIFICADO"),
    @("CC", 1048457729, "DAYANA", "URRUCHURTU NIÑO", "CERTIFICADO"),
    @("TI", 1049019898, "KAREN YURLEIDY", "MARIN VARGAS", "RETIRO VOLUNTARIO"),
    @("CC", 1087985197, "GISELL MARIANA", "MARIN LARROTA", "CERTIFICADO"),
    @("CC", 1096184002, "DANIELA", "ROJAS BOTELLO", "CERTIFICADO"),
    @("CC", 1096186262, "KEVIN ANDRES", "PARADA SUAREZ", "RETIRO VOLUNTARIO"),
    @("CC", 1096189477, "KAMILA", "QUINTERO CARREÑO", "CERTIFICADO"),
    @("CC", 1097183074, "MARIA JOSE", "ORTIZ GUIZA", "CERTIFICADO"),
    @("CC", 1144182405, "CAROLAIN", "ABANIS PEREZ", "CERTIFICADO"),
    @("CC", 63469380, "VIDA EMPERATRIZ", "SANTOS YAIN", "CERTIFICADO")
)

$r = 16
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
